$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 44523.84451388889
$ws.Range("C6").Value = 44523.84770833334
$ws.Range("D6").Value = "IP Address"
$ws.Range("E6").Value = 100
$ws.Range("F6").Value = 275
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = 44523.84771990741
$ws.Range("I6").Value = "1pogus"
$ws.Range("J6").Value = "library(readxl)
cmv <- data/cmv.xlsx
cmv %>%
  cmv_tidy <- pivot_longer(donor_negative:donor_positive, names_to = ""donor_status"")"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 44523.84465277778
$ws.Range("C7").Value = 44523.8491087963
$ws.Range("D7").Value = "IP Address"
$ws.Range("E7").Value = 100
$ws.Range("F7").Value = 385
$ws.Range("G7").Value = $true
$ws.Range("H7").Value = 44523.84912037037
$ws.Range("I7").Value = "3mihar"
$ws.Range("J7").Value = "cmv <- read_excel(""data/cmv.xlsx"")
cmv_subset <- cmv %>% 
filter(age > 65)
write_csv(cmv_subset, ""data/cmv_subset.csv"")
cmv_tidy <- cmv %>% 
pivot_longer(cols = donor_negative:donor_positive, names_to = ""donor_status"", values_to = ""recipient_status"") %>% 
drop_na()
cmv_tidy %>% 
group_by(cmv) %>% 
summarize(mean_age = mean(age))"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 44523.84481481482
$ws.Range("C8").Value = 44523.85103009259
$ws.Range("D8").Value = "IP Address"
$ws.Range("E8").Value = 100
$ws.Range("F8").Value = 537
$ws.Range("G8").Value = $true
$ws.Range("H8").Value = 44523.85103009259
$ws.Range("I8").Value = "2dunic"
$ws.Range("J8").Value = "library(readxl)
cmv <- read.table(data/cmv.xlsx)
cmv_subset <- cmv(col = 'age' = >64)
cmv_tidy <- cmv_subset %>%
pivot_longer('0':last_col(), names_to = ""donor_status"", ""recipient_status"")
cmv_tidy
avarage(age)"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 44523.84476851852
$ws.Range("C9").Value = 44523.85337962963
$ws.Range("D9").Value = "IP Address"
$ws.Range("E9").Value = 100
$ws.Range("F9").Value = 743
$ws.Range("G9").Value = $true
$ws.Range("H9").Value = 44523.85337962963
$ws.Range("I9").Value = "1dabec"
$ws.Range("J9").Value = "library(readxl)
cmv <- read_xlsx(""data/cmv.xlsx"")
cmv_subset <- cmv %>% filter(age > 65)
cmv %>%
 pivot_longer(``donor_negative``:last_col(), names_to = ""donor_status"", values_to = ""recipient_status"") %>%
 drop_na()
grouped_data <- group_by(cmv_tidy, cmv)
summarise(grouped_data, average = mean(age))"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 44523.84472222222
$ws.Range("C10").Value = 44523.86057870371
$ws.Range("D10").Value = "IP Address"
$ws.Range("E10").Value = 100
$ws.Range("F10").Value = 1370
$ws.Range("G10").Value = $true
$ws.Range("H10").Value = 44523.86059027778
$ws.Range("I10").Value = "1bimil"
$ws.Range("J10").Value = "cmv <- read_excel(""data/cmv.xlsx"")
cmv_subset <- cmv %>% 
  filter(age > 65) %>% 
  write_csv(""data/cmv_subset.csv"")
cmv_tidy <- cmv %>% 
  pivot_longer(donor_negative:donor_positive, names_to = ""donor_status"", values_to = ""recipient_status"") %>% 
  drop_na()
cmv_tidy %>% 
  group_by(cmv) %>% 
  summarise(average_age = mean(age))
"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 44523.86260416667
$ws.Range("C11").Value = 44523.86299768518
$ws.Range("D11").Value = "IP Address"
$ws.Range("E11").Value = 100
$ws.Range("F11").Value = 34
$ws.Range("G11").Value = $true
$ws.Range("H11").Value = 44523.86300925926
$ws.Range("I11").Value = "1davec"
$ws.Range("J11").Value = "
cmv <- read_excel(""data/cmv.xlsx"")
cmv_subset <- cmv %>% filter(age > 65)
write_csv(cmv_subset, ""data/cmv_subset.csv"")
cmv <- read_excel(""cmv.xlsx"")
cmv_tidy <- cmv %>% pivot_longer(``donor_negative``: last_col()) %>% separate(value, into = c(""donor_status"", ""recipient_status""), sep= 0) %>% drop_na()
cmv_tidy <- cmv_tidy[-c(7)]
cmv_tidy <- cmv_tidy %>% rename(donor_status = name)"

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 44518.82976851852
$ws.Range("C12").Value = 44518.82983796296
$ws.Range("D12").Value = "IP Address"
$ws.Range("E12").Value = 50
$ws.Range("F12").Value = 6
$ws.Range("G12").Value = $false
$ws.Range("H12").Value = 44523.87890046297
$ws.Range("I12").Value = "3kusou"

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 44518.84795138889
$ws.Range("C13").Value = 44518.85869212963
$ws.Range("D13").Value = "IP Address"
$ws.Range("E13").Value = 50
$ws.Range("F13").Value = 928
$ws.Range("G13").Value = $false
$ws.Range("H13").Value = 44523.87898148148
$ws.Range("I13").Value = "0garbc"

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 44523.84482638889
$ws.Range("C14").Value = 44523.84958333333
$ws.Range("D14").Value = "IP Address"
$ws.Range("E14").Value = 50
$ws.Range("F14").Value = 411
$ws.Range("G14").Value = $false
$ws.Range("H14").Value = 44523.87898148148
$ws.Range("I14").Value = "3+ka+mei"

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 44523.18445601852
$ws.Range("C15").Value = 44523.29313657407
$ws.Range("D15").Value = "Spam"
$ws.Range("E15").Value = 50
$ws.Range("F15").Value = 9390
$ws.Range("G15").Value = $false
$ws.Range("H15").Value = 44523.87903935185
$ws.Range("I15").Value = "2nesch"

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 44523.84469907408
$ws.Range("C16").Value = 44523.85349537037
$ws.Range("D16").Value = "IP Address"
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 759
$ws.Range("G16").Value = $false
$ws.Range("H16").Value = 44523.87905092593
$ws.Range("I16").Value = "3hostc"
